$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab name) to fix the typo "Cohort_Retentio" -> "Cohort_Retention"
$ws.Name = "Cohort_Retention"

# Update the active cell selection on the sheet from G9 to K7
$ws.Range("K7").Select()
